$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.43
$ws.Range("E2").Value = 1.43

$ws.Range("B3").Value = 1.4
$ws.Range("E3").Value = 1.31
$ws.Range("F3").Value = 1.29

$ws.Range("C4").Value = 1.45
$ws.Range("D4").Value = 1.36
$ws.Range("F4").Value = 0.98

$ws.Range("B5").Value = 1.36
$ws.Range("C5").Value = 1.34
$ws.Range("G5").Value = 0.53

$ws.Range("C6").Value = 1.38
$ws.Range("D6").Value = 1.64
$ws.Range("G6").Value = 1.06

$ws.Range("F7").Value = 1.47
